# Auto-generated edit script applying numeric corrections to the Odin_Profits data
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per the commit diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2408.5454
$ws.Range("I9").Value = 2956.1428
$ws.Range("J9").Value = 1450.25
$ws.Range("K9").Value = 2956.1428
$ws.Range("L9").Value = 1450.25
$ws.Range("M9").Value = -2787.1428
$ws.Range("N9").Value = -1788.25
$ws.Range("H15").Value = 169749.95
$ws.Range("I15").Value = 169749.95
$ws.Range("K15").Value = 509249.85
$ws.Range("M15").Value = -509080.85
$ws.Range("H33").Value = 698.8889
$ws.Range("I33").Value = 612.8570999999999
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 612.8570999999999
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -383.8570999999999
$ws.Range("N33").Value = -1458
$ws.Range("H64").Value = 27781606
$ws.Range("I64").Value = 55557616
$ws.Range("K64").Value = 55557616
$ws.Range("M64").Value = -55557368
$ws.Range("H67").Value = 27781606
$ws.Range("I67").Value = 55557616
$ws.Range("K67").Value = 55557616
$ws.Range("M67").Value = -55556758
$ws.Range("H98").Value = 17994
$ws.Range("I98").Value = 17994
$ws.Range("K98").Value = 17994
$ws.Range("M98").Value = -16496
$ws.Range("H112").Value = 1991.7675
$ws.Range("J112").Value = 2128.8948
$ws.Range("L112").Value = 6386.6844
$ws.Range("N112").Value = -8602.6844
$ws.Range("H122").Value = 17994
$ws.Range("I122").Value = 17994
$ws.Range("K122").Value = 53982
$ws.Range("M122").Value = -51532
$ws.Range("H132").Value = 261379.97
$ws.Range("I132").Value = 281345.7
$ws.Range("J132").Value = 16799.75
$ws.Range("K132").Value = 844037.1000000001
$ws.Range("L132").Value = 50399.25
$ws.Range("M132").Value = -841507.1000000001
$ws.Range("N132").Value = -55459.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1606.25
$ws.Range("I45").Value = 1370
$ws.Range("K45").Value = 1370
$ws.Range("M45").Value = -993
$ws.Range("H61").Value = 4487.1724
$ws.Range("I61").Value = 4436.684
$ws.Range("K61").Value = 4436.684
$ws.Range("M61").Value = -4224.684
$ws.Range("H63").Value = 5224.375
$ws.Range("J63").Value = 5455.1177
$ws.Range("L63").Value = 5455.1177
$ws.Range("N63").Value = -6827.1177
$ws.Range("H66").Value = 5224.375
$ws.Range("J66").Value = 5455.1177
$ws.Range("L66").Value = 27275.5885
$ws.Range("N66").Value = -34139.5885
$ws.Range("H74").Value = 6215.3076
$ws.Range("I74").Value = 6579.9
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 6579.9
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -5705.9
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 6215.3076
$ws.Range("I77").Value = 6579.9
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 32899.5
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -28531.5
$ws.Range("N77").Value = -33736
$ws.Range("H122").Value = 2959.2
$ws.Range("I122").Value = 2662.2666
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 7986.7998
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -5536.7998
$ws.Range("N122").Value = -16450
$ws.Range("H132").Value = 696958.4399999999
$ws.Range("I132").Value = 721086.0600000001
$ws.Range("K132").Value = 2163258.18
$ws.Range("M132").Value = -2160728.18
$ws.Range("H136").Value = 4487.1724
$ws.Range("I136").Value = 4436.684
$ws.Range("K136").Value = 13310.052
$ws.Range("M136").Value = -10760.052

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4703.885
$ws.Range("I94").Value = 1775
$ws.Range("J94").Value = 7214.357
$ws.Range("K94").Value = 1775
$ws.Range("L94").Value = 7214.357
$ws.Range("M94").Value = -1324
$ws.Range("N94").Value = -8116.357
$ws.Range("H99").Value = 7248.353
$ws.Range("I99").Value = 6786.646
$ws.Range("K99").Value = 6786.646
$ws.Range("M99").Value = -5288.646
$ws.Range("H105").Value = 4995.4287
$ws.Range("I105").Value = 4994.6665
$ws.Range("K105").Value = 4994.6665
$ws.Range("M105").Value = -3247.6665
$ws.Range("H134").Value = 930543.9
$ws.Range("I134").Value = 1045535.25
$ws.Range("K134").Value = 3136605.75
$ws.Range("M134").Value = -3134070.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23821290
$ws.Range("I58").Value = 37045136
$ws.Range("K58").Value = 37045136
$ws.Range("M58").Value = -37044933
$ws.Range("H62").Value = 10195.111
$ws.Range("J62").Value = 6609.6665
$ws.Range("L62").Value = 6609.6665
$ws.Range("N62").Value = -7857.6665
$ws.Range("H65").Value = 10195.111
$ws.Range("J65").Value = 6609.6665
$ws.Range("L65").Value = 33048.3325
$ws.Range("N65").Value = -39288.3325
$ws.Range("H136").Value = 23821290
$ws.Range("I136").Value = 37045136
$ws.Range("K136").Value = 111135408
$ws.Range("M136").Value = -111132858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2403.4211
$ws.Range("I132").Value = 1100
$ws.Range("J132").Value = 3351.3635
$ws.Range("K132").Value = 9900
$ws.Range("L132").Value = 30162.2715
$ws.Range("M132").Value = -7370
$ws.Range("N132").Value = -35222.2715
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H136").Value = 23814538
$ws.Range("I136").Value = 23814538
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 71443614
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -71438514
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 20835994
$ws.Range("J139").Value = 5949
$ws.Range("L139").Value = 17847
$ws.Range("N139").Value = -28127

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 30001
$ws.Range("J49").Value = 30001
$ws.Range("L49").Value = 30001
$ws.Range("N49").Value = -30369
$ws.Range("H132").Value = 10511
$ws.Range("I132").Value = 11104.75
$ws.Range("J132").Value = 6948.5
$ws.Range("K132").Value = 33314.25
$ws.Range("L132").Value = 20845.5
$ws.Range("M132").Value = -30784.25
$ws.Range("N132").Value = -25905.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9990.096
$ws.Range("I7").Value = 7600.077
$ws.Range("J7").Value = 13873.875
$ws.Range("K7").Value = 7600.077
$ws.Range("L7").Value = 13873.875
$ws.Range("M7").Value = -7488.077
$ws.Range("N7").Value = -14097.875
$ws.Range("H40").Value = 40441.89
$ws.Range("I40").Value = 111332.664
$ws.Range("J40").Value = 4996.5
$ws.Range("K40").Value = 111332.664
$ws.Range("L40").Value = 4996.5
$ws.Range("M40").Value = -111196.664
$ws.Range("N40").Value = -5268.5
$ws.Range("H100").Value = 3138.52
$ws.Range("I100").Value = 4278
$ws.Range("J100").Value = 2378.8667
$ws.Range("K100").Value = 4278
$ws.Range("L100").Value = 2378.8667
$ws.Range("M100").Value = -3737
$ws.Range("N100").Value = -3460.8667
$ws.Range("H126").Value = 9990.096
$ws.Range("I126").Value = 7600.077
$ws.Range("J126").Value = 13873.875
$ws.Range("K126").Value = 22800.231
$ws.Range("L126").Value = 41621.625
$ws.Range("M126").Value = -20330.231
$ws.Range("N126").Value = -46561.625
$ws.Range("H132").Value = 2351
$ws.Range("I132").Value = 1910.4
$ws.Range("K132").Value = 5731.200000000001
$ws.Range("M132").Value = -3201.200000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5993.2915
$ws.Range("I62").Value = 7970.4
$ws.Range("J62").Value = 5473
$ws.Range("K62").Value = 7970.4
$ws.Range("L62").Value = 5473
$ws.Range("M62").Value = -7346.4
$ws.Range("N62").Value = -6721
$ws.Range("H65").Value = 5993.2915
$ws.Range("I65").Value = 7970.4
$ws.Range("J65").Value = 5473
$ws.Range("K65").Value = 39852
$ws.Range("L65").Value = 27365
$ws.Range("M65").Value = -36732
$ws.Range("N65").Value = -33605
$ws.Range("H132").Value = 10049.889
$ws.Range("I132").Value = 5981.593
$ws.Range("J132").Value = 22254.777
$ws.Range("K132").Value = 17944.779
$ws.Range("L132").Value = 66764.33099999999
$ws.Range("M132").Value = -15414.779
$ws.Range("N132").Value = -71824.33099999999
$ws.Range("H136").Value = 15636282
$ws.Range("I136").Value = 18527534
$ws.Range("J136").Value = 23517
$ws.Range("K136").Value = 55582602
$ws.Range("L136").Value = 70551
$ws.Range("M136").Value = -55580052
$ws.Range("N136").Value = -75651

Write-Host "Applied 218 cell updates across 8 sheets"
